$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "606.16"); set it via
# a temporary text NumberFormat so the COM layer keeps it as a text value (matching
# the source inlineStr cells) and then clear the format again so no style lingers.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '65.761.97'
$ws.Range("E2").Value = '  +1.00%  '
Set-TextValue $ws.Range("D3") '2.701.36'
$ws.Range("E3").Value = '  +1.73%  '
$ws.Range("E4").Value = '  +0.03%  '
Set-TextValue $ws.Range("D5") '606.16'
$ws.Range("E5").Value = '  +1.95%  '
Set-TextValue $ws.Range("D6") '157.59'
$ws.Range("E6").Value = '  +0.89%  '
$ws.Range("E7").Value = '  +0.08%  '
Set-TextValue $ws.Range("D8") '0.587'
$ws.Range("E8").Value = '  -1.13%  '
$ws.Range("E9").Value = '  +4.99%  '
Set-TextValue $ws.Range("D10") '6.05'
Set-TextValue $ws.Range("D11") '0.401'
$ws.Range("E11").Value = '  +0.33%  '
$ws.Range("E12").Value = '  +1.06%  '
Set-TextValue $ws.Range("D13") '30.07'
$ws.Range("E13").Value = '  +3.56%  '
Set-TextValue $ws.Range("D14") '0.0000204'
$ws.Range("E14").Value = '  +8.86%  '
Set-TextValue $ws.Range("D15") '3.186.34'
$ws.Range("E15").Value = '  +1.75%  '
Set-TextValue $ws.Range("D16") '65.664.77'
$ws.Range("E16").Value = '  +1.00%  '
Set-TextValue $ws.Range("D17") '2.718.41'
$ws.Range("E17").Value = '  +3.54%  '
Set-TextValue $ws.Range("D18") '12.71'
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("E19").Value = '  +1.25%  '
Set-TextValue $ws.Range("D20") '359.89'
$ws.Range("E20").Value = '  +1.57%  '
Set-TextValue $ws.Range("D21") '7.51'
$ws.Range("E21").Value = '  +3.11%  '
$ws.Range("E22").Value = '  -0.09%  '
Set-TextValue $ws.Range("D23") '70.36'
$ws.Range("E23").Value = '  +3.13%  '
Set-TextValue $ws.Range("D24") '9.84'
$ws.Range("E24").Value = '  +3.27%  '
$ws.Range("E25").Value = '  +11.56%  '
$ws.Range("E26").Value = '  -4.60%  '
Set-TextValue $ws.Range("D27") '1.69'
$ws.Range("E27").Value = '  +2.97%  '
$ws.Range("E28").Value = '  +3.76%  '
Set-TextValue $ws.Range("D29") '8.32'
$ws.Range("E29").Value = '  +1.77%  '
$ws.Range("E30").Value = '  +4.11%  '
$ws.Range("E31").Value = '  +0.06%  '
Set-TextValue $ws.Range("D32") '540.44'
$ws.Range("E32").Value = '  +3.18%  '
$ws.Range("E33").Value = '  +0.05%  '
Set-TextValue $ws.Range("D34") '6.67'
$ws.Range("E34").Value = '  +4.59%  '
Set-TextValue $ws.Range("D35") '5.40'
$ws.Range("E35").Value = '  -4.45%  '
$ws.Range("E36").Value = '  +1.01%  '
Set-TextValue $ws.Range("D37") '20.72'
$ws.Range("E37").Value = '  +2.19%  '
Set-TextValue $ws.Range("D38") '162.59'
$ws.Range("E39").Value = '  -1.32%  '
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("E42").Value = '  +1.47%  '
Set-TextValue $ws.Range("D43") '167.96'
$ws.Range("E43").Value = '  +1.47%  '
Set-TextValue $ws.Range("D44") '4.19'
$ws.Range("E44").Value = '  +2.03%  '
Set-TextValue $ws.Range("D45") '0.0616'
$ws.Range("E45").Value = '  -0.22%  '
Set-TextValue $ws.Range("D46") '23.61'
$ws.Range("E46").Value = '  +2.76%  '
$ws.Range("E47").Value = '  +2.35%  '
$ws.Range("E48").Value = '  +4.32%  '
Set-TextValue $ws.Range("D49") '0.659'
$ws.Range("E49").Value = '  +1.44%  '
Set-TextValue $ws.Range("D50") '21.05'
$ws.Range("E50").Value = '  +7.68%  '
Set-TextValue $ws.Range("D51") '0.0984'
$ws.Range("E51").Value = '  -0.36%  '
